$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix A79: the recorded timestamp was wrong; correct it (R script re-run) ---
$ws.Range("A79").Value = 45450.2916666667

# --- Append new row 80 with the latest scraped OHLCV data ---
$ws.Range("A80").Value = 45453.3480208333
$ws.Range("B80").Value = 2100
$ws.Range("C80").Value = 6.38000011444092
$ws.Range("D80").Value = 6.26000022888184
$ws.Range("E80").Value = 6.38000011444092
$ws.Range("F80").Value = 6.26000022888184

# G80 ("adj_close") is stored as text in this sheet (matches existing column
# G cells, which are all shared-string numbers) - use a quote-prefix so it
# is written out as a literal string value instead of a number.
$ws.Range("G80").Value = "'6.26000022888184"

# H80 ("ticker")
$ws.Range("H80").Value = "PAL.MI"

# Give A80 the same date/time number format used by the rest of column A
$ws.Range("A79").Copy()
$ws.Range("A80").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Column G cells elsewhere carry the default (General) style - undo the
# quote-prefix style bump on G80 so it matches its neighbours exactly.
$ws.Range("G80").Style = "Normal"
